$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "515.10", "0.500") are preserved exactly as text, matching the source data.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "57.876.29"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "3.057.62"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "515.10"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "141.25"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.435"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "0.375"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "3.583.19"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "26.24"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").Value = "0.0000163"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "57.899.63"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "3.061.89"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "6.09"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "8.14"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "331.08"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "0.500"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "65.43"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "0.0₃0901"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("D28").Value = "6.45"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").Value = "  +5.96%  "
$ws.Range("D30").Value = "1.80"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").Value = "20.57"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "154.49"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "4.51"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "5.99"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").Value = "27.08"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").Value = "1.26"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("D38").Value = "0.0676"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").Value = "3.100.53"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("D40").Value = "3.90"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").Value = "36.66"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "0.656"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "2.288.62"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("D45").Value = "0.0256"
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("D46").Value = "1.37"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "20.65"
$ws.Range("E47").Value = "  +4.29%  "
$ws.Range("D48").Value = "0.939"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "5.92"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").Value = "0.729"
$ws.Range("E50").Value = "  +8.20%  "
$ws.Range("D51").Value = "0.0875"
$ws.Range("E51").Value = "  +2.01%  "

# Restore default cell style (removes the explicit text-number-format
# so the cells match the original unstyled appearance).
$priceRange.Style = "Normal"

